$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Lugo" and "Almeria" rows (row 47 and row 48).
# Row 47 was Lugo (58,5,53,1); Row 48 was Almeria (58,72,53,1).
# After the edit, row 47 should be Almeria (58,72,53,1) and row 48 Lugo (58,5,53,1).
$ws.Range("A47").Value = "Almeria"
$ws.Range("C47").Value = 72
$ws.Range("A48").Value = "Lugo"
$ws.Range("C48").Value = 5

# Update the "last updated" timestamp string in cell A1.
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 20:16"
